$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = "CreateAccount"
$ws.Range("C6").Value = "ViewMessage"
$ws.Range("B8").Value = "CreateGroup"
$ws.Range("C8").Value = "JoinGroup"

$ws.Range("C8").Select()
